# Add new rows of travel notes to the "買" (shopping) and "吃" (food) sheets,
# then switch the active tab to "吃".

$wb = $excel.ActiveWorkbook

# ---- Sheet "買": two new rows (16-17) ----
$wsBuy = $wb.Worksheets.Item("買")

# Shared-string creation order matters (mirrors how the original author typed
# the cells), so write C17, then C16, then A16.
$wsBuy.Range("C17").Value = "B1八重洲地下中央口的東京車站一番街"
$wsBuy.Range("C16").Value = "東京車站丸之內北口對面"
$wsBuy.Range("A16").Value = "東京車站KITTE"

# ---- Sheet "吃": five new rows (5-9) ----
$wsEat = $wb.Worksheets.Item("吃")

$wsEat.Range("A5").Value = "川上庵"
$wsEat.Range("D5").Value = "長野県北佐久郡軽井沢町軽井沢6-10"
$wsEat.Range("C5").Value = "11:00 AM～10:00 PM"

$wsEat.Range("A6").Value = "村民食堂 [長倉]"
$wsEat.Range("C6").Value = "11:00 AM～10:00 PM"

$wsEat.Range("A7").Value = "明治亭"
$wsEat.Range("D7").Value = "輕井澤王子購物廣場味之街 AJ-12‧13"

$wsEat.Range("A8").Value = "Aging Beef 熟成和牛"

$wsEat.Range("A9").Value = "御曹司きよやす庵"
$wsEat.Range("D9").Value = "長野県北佐久郡軽井沢町軽井沢1178-161"

# New font for D5 (grey "微軟正黑體" note text)
$wsEat.Range("D5").Font.Name = "微軟正黑體"
$wsEat.Range("D5").Font.Size = 12
$wsEat.Range("D5").Font.Color = 9276813

# ---- View state: sheet "吃" becomes the active/visible tab ----
$wsBuy.Range("A16").Select()
$wb.Worksheets.Item("吃").Activate()
$wsEat.Range("D9").Select()
